# Update latest output (run 75)
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Schedule": refresh the single pump-run row and append the next run
# ---------------------------------------------------------------------------
$schedule = $wb.Worksheets.Item("Schedule")

# Row 2 gets new values (the optimisation re-ran and produced different numbers)
$schedule.Cells.Item(2, 1).Value2 = 46040.29166666666
$schedule.Cells.Item(2, 2).Value2 = 46040.79166666666
$schedule.Cells.Item(2, 3).Value2 = 12
$schedule.Cells.Item(2, 4).Value2 = 45.36
$schedule.Cells.Item(2, 5).Value2 = -59.48107424999998
$schedule.Cells.Item(2, 6).Value2 = -1.311311160714285

# Row 3 is a brand new pump run, appended below the existing one
$schedule.Cells.Item(3, 1).Value2 = 46040.83333333334
$schedule.Cells.Item(3, 2).Value2 = 46041
$schedule.Cells.Item(3, 3).Value2 = 4
$schedule.Cells.Item(3, 4).Value2 = 15.12
$schedule.Cells.Item(3, 5).Value2 = 370.4505674999999
$schedule.Cells.Item(3, 6).Value2 = 24.50069890873015

# Match the date/time formatting used by row 2 for the new Start/Stop Time cells
$schedule.Range("A3:B3").NumberFormat = $schedule.Range("A2:B2").NumberFormat

# ---------------------------------------------------------------------------
# Sheet "Detailed": updated forecast/historical price curve & pump status
# ---------------------------------------------------------------------------
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Cells.Item(9, 5).Value = "OFF"
$detailed.Cells.Item(10, 5).Value = "OFF"

$detailed.Cells.Item(11, 2).Value2 = 57.06003
$detailed.Cells.Item(11, 5).Value = "OFF"

$detailed.Cells.Item(12, 2).Value2 = 57.06003
$detailed.Cells.Item(12, 5).Value = "OFF"

$detailed.Cells.Item(13, 2).Value2 = 55.33101
$detailed.Cells.Item(13, 3).Value = "historical"
$detailed.Cells.Item(13, 5).Value = "OFF"

$detailed.Cells.Item(14, 2).Value2 = 56.26125
$detailed.Cells.Item(14, 3).Value = "historical"
$detailed.Cells.Item(14, 5).Value = "OFF"

$detailed.Cells.Item(15, 5).Value = "OFF"

$detailed.Cells.Item(16, 2).Value2 = 36.06
$detailed.Cells.Item(17, 2).Value2 = 0.61797
$detailed.Cells.Item(18, 2).Value2 = -5.4274
$detailed.Cells.Item(19, 2).Value2 = 0
$detailed.Cells.Item(20, 2).Value2 = -4.66577
$detailed.Cells.Item(21, 2).Value2 = -5.31676
$detailed.Cells.Item(22, 2).Value2 = -3.6481
$detailed.Cells.Item(23, 2).Value2 = 0.64369
$detailed.Cells.Item(24, 2).Value2 = 0.00976
$detailed.Cells.Item(25, 2).Value2 = -4.66482
$detailed.Cells.Item(26, 2).Value2 = 0.00025
$detailed.Cells.Item(27, 2).Value2 = -5.51
$detailed.Cells.Item(28, 2).Value2 = -5.36283
$detailed.Cells.Item(29, 2).Value2 = -6.0715
$detailed.Cells.Item(30, 2).Value2 = -13.5
$detailed.Cells.Item(31, 2).Value2 = -20.57961
$detailed.Cells.Item(32, 2).Value2 = -13.5
$detailed.Cells.Item(33, 2).Value2 = -11.52252
$detailed.Cells.Item(34, 2).Value2 = -7.11948
$detailed.Cells.Item(35, 2).Value2 = -6.60693
$detailed.Cells.Item(37, 2).Value2 = 0.31443
$detailed.Cells.Item(38, 2).Value2 = 3.98303
$detailed.Cells.Item(39, 2).Value2 = 16.86036

$detailed.Cells.Item(40, 2).Value2 = 40.99071
$detailed.Cells.Item(40, 5).Value = "OFF"

$detailed.Cells.Item(41, 2).Value2 = 53.66969

$detailed.Cells.Item(42, 5).Value = "ON"

$detailed.Cells.Item(43, 2).Value2 = 57.3
$detailed.Cells.Item(43, 5).Value = "ON"

$detailed.Cells.Item(44, 2).Value2 = 56.94617
$detailed.Cells.Item(44, 5).Value = "ON"

$detailed.Cells.Item(45, 2).Value2 = 43.56065
$detailed.Cells.Item(45, 5).Value = "ON"

$detailed.Cells.Item(46, 2).Value2 = 56.66228
$detailed.Cells.Item(46, 5).Value = "ON"

$detailed.Cells.Item(47, 5).Value = "ON"
$detailed.Cells.Item(48, 5).Value = "ON"
$detailed.Cells.Item(49, 5).Value = "ON"
